# Update the "Exchange Offer" discount amount from 40000 to 35000
# for every data row (rows 3-36) on the "Sheet1" worksheet, column O.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("O3:O36").Value = 35000
